# Weekly fruit/vegetable price update: insert 3 new sample rows at the top
# of the Espárragos block (pushing the existing 147:181 rows down to
# 150:184) and populate them with the latest observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 147; everything below shifts
# down by three (old row 147 -> new row 150, ..., old row 181 -> new row 184).
$ws.Rows("147:149").Insert()

# New row 147: Banquete, $/bandeja 10 kilos
$ws.Range("A147").Value = 9
$ws.Range("B147").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C147").Value = "Metropolitana"
$ws.Range("D147").Value = 45218
$ws.Range("E147").Value = 13
$ws.Range("F147").Value = 300000000
$ws.Range("G147").Value = "Espárragos"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Banquete"
$ws.Range("J147").Value = 43
$ws.Range("K147").Value = 15000
$ws.Range("L147").Value = 15000
$ws.Range("M147").Value = 15000
$ws.Range("N147").Value = "$/bandeja 10 kilos"
$ws.Range("O147").Value = "Provincia de Linares"
$ws.Range("P147").Value = 1500
$ws.Range("Q147").Value = 10
$ws.Range("R147").Value = "Hortaliza"

# New row 148: Primera, $/bandeja 10 kilos
$ws.Range("A148").Value = 9
$ws.Range("B148").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C148").Value = "Metropolitana"
$ws.Range("D148").Value = 45218
$ws.Range("E148").Value = 13
$ws.Range("F148").Value = 300000000
$ws.Range("G148").Value = "Espárragos"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 70
$ws.Range("K148").Value = 13000
$ws.Range("L148").Value = 13000
$ws.Range("M148").Value = 13000
$ws.Range("N148").Value = "$/bandeja 10 kilos"
$ws.Range("O148").Value = "Provincia de Linares"
$ws.Range("P148").Value = 1300
$ws.Range("Q148").Value = 10
$ws.Range("R148").Value = "Hortaliza"

# New row 149: Segunda, $/bandeja 10 kilos
$ws.Range("A149").Value = 9
$ws.Range("B149").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C149").Value = "Metropolitana"
$ws.Range("D149").Value = 45218
$ws.Range("E149").Value = 13
$ws.Range("F149").Value = 300000000
$ws.Range("G149").Value = "Espárragos"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Segunda"
$ws.Range("J149").Value = 52
$ws.Range("K149").Value = 11000
$ws.Range("L149").Value = 11000
$ws.Range("M149").Value = 11000
$ws.Range("N149").Value = "$/bandeja 10 kilos"
$ws.Range("O149").Value = "Provincia de Linares"
$ws.Range("P149").Value = 1100
$ws.Range("Q149").Value = 10
$ws.Range("R149").Value = "Hortaliza"
